# Fix creation of requirements: clear the sample/template requirement rows
# (CODIGO / DESCRIPCION / STORY POINTS / horas values) from the first sheet
# so that only the header row remains, and reset the sheet's selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the sample data that used to populate rows 2-18 (columns A-D),
# leaving just the header row (row 1) and the blank template rows below.
$ws.Range("A2:D18").ClearContents()

# Reset the active selection, as in the updated template.
$ws.Range("A9").Select()
